$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) and Volume(1h) (column E) cells per the refreshed crypto snapshot.
# Cells whose new text would otherwise be auto-parsed by Excel as a number (e.g. "1.003")
# are explicitly kept as Text so the stored value remains the literal string.
$ws.Cells.Item(2, 4).Value = "25.745.73"
$ws.Cells.Item(2, 5).Value = "  -1.96%  "
$ws.Cells.Item(3, 4).Value = "1.611.71"
$ws.Cells.Item(3, 5).Value = "  -3.90%  "
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.003"
$ws.Cells.Item(4, 5).Value = "  +0.07%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "207.80"
$ws.Cells.Item(5, 5).Value = "  -1.98%  "
$ws.Cells.Item(6, 5).Value = "  -1.51%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "1.005"
$ws.Cells.Item(7, 5).Value = "  +0.23%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.2548"
$ws.Cells.Item(8, 5).Value = "  -4.16%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.06184"
$ws.Cells.Item(9, 5).Value = "  -1.80%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "20.03"
$ws.Cells.Item(10, 5).Value = "  -6.46%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.07513"
$ws.Cells.Item(11, 5).Value = "  -0.56%  "
$ws.Cells.Item(12, 4).Value = "1.615.92"
$ws.Cells.Item(12, 5).Value = "  -3.69%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "4.345"
$ws.Cells.Item(13, 5).Value = "  -2.83%  "
$ws.Cells.Item(14, 4).Value = "1.842.56"
$ws.Cells.Item(14, 5).Value = "  -3.34%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.5388"
$ws.Cells.Item(15, 5).Value = "  -4.40%  "
$ws.Cells.Item(16, 4).Value = "0.0₅7788"
$ws.Cells.Item(16, 5).Value = "  -3.21%  "
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "63.87"
$ws.Cells.Item(17, 5).Value = "  -4.79%  "
$ws.Cells.Item(18, 4).Value = "25.735.88"
$ws.Cells.Item(18, 5).Value = "  -1.11%  "
$ws.Cells.Item(19, 5).Value = "  +0.16%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "4.572"
$ws.Cells.Item(20, 5).Value = "  -5.34%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "182.69"
$ws.Cells.Item(21, 5).Value = "  -3.06%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "9.956"
$ws.Cells.Item(22, 5).Value = "  -4.65%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "1.006"
$ws.Cells.Item(23, 5).Value = "  +0.15%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "6.000"
$ws.Cells.Item(24, 5).Value = "  -3.22%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "144.17"
$ws.Cells.Item(25, 5).Value = "  -3.98%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "0.1194"
$ws.Cells.Item(26, 5).Value = "  -5.03%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "7.293"
$ws.Cells.Item(27, 5).Value = "  -4.06%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "15.36"
$ws.Cells.Item(28, 5).Value = "  -4.25%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "1.356"
$ws.Cells.Item(29, 5).Value = "  -0.35%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "0.05876"
$ws.Cells.Item(30, 5).Value = "  -5.57%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "1.233"
$ws.Cells.Item(31, 5).Value = "  -4.16%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "3.337"
$ws.Cells.Item(32, 5).Value = "  -4.94%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "3.305"
$ws.Cells.Item(33, 5).Value = "  -4.20%  "
$ws.Cells.Item(34, 5).Value = "  -3.40%  "
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "0.9569"
$ws.Cells.Item(35, 5).Value = "  -4.72%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "2.386"
$ws.Cells.Item(36, 5).Value = "  -0.85%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "2.701"
$ws.Cells.Item(37, 5).Value = "  -1.33%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.5650"
$ws.Cells.Item(38, 5).Value = "  -7.04%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.01572"
$ws.Cells.Item(39, 5).Value = "  -3.14%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "1.003"
$ws.Cells.Item(40, 5).Value = "  -0.33%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.8288"
$ws.Cells.Item(41, 5).Value = "  -4.94%  "
$ws.Cells.Item(42, 4).Value = "1.013.75"
$ws.Cells.Item(42, 5).Value = "  -8.18%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "5.572"
$ws.Cells.Item(43, 5).Value = "  -8.73%  "
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "98.74"
$ws.Cells.Item(44, 5).Value = "  -1.33%  "
$ws.Cells.Item(45, 4).Value = "1.766.02"
$ws.Cells.Item(45, 5).Value = "  -3.32%  "
$ws.Cells.Item(46, 5).Value = "  -2.07%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "1.002"
$ws.Cells.Item(47, 5).Value = "  -0.23%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "53.70"
$ws.Cells.Item(48, 5).Value = "  -4.41%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "0.05158"
$ws.Cells.Item(49, 5).Value = "  -1.49%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "7.812"
$ws.Cells.Item(50, 5).Value = "  -3.01%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.4216"
$ws.Cells.Item(51, 5).Value = "  -0.92%  "
